$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 225.86667
$ws.Range("I11").Value = 225.86667
$ws.Range("K11").Value = 225.86667
$ws.Range("M11").Value = -85.86667
$ws.Range("H17").Value = 5031.0645
$ws.Range("J17").Value = 5031.0645
$ws.Range("L17").Value = 15093.1935
$ws.Range("N17").Value = -15429.1935
$ws.Range("H33").Value = 364.23077
$ws.Range("I33").Value = 359.77274
$ws.Range("K33").Value = 359.77274
$ws.Range("M33").Value = -130.77274
$ws.Range("H111").Value = 719.46155
$ws.Range("I111").Value = 586.63635
$ws.Range("K111").Value = 1759.90905
$ws.Range("M111").Value = 1307.09095
$ws.Range("H112").Value = 2043.2222
$ws.Range("J112").Value = 2301.2856
$ws.Range("L112").Value = 6903.8568
$ws.Range("N112").Value = -9119.856800000001
$ws.Range("H116").Value = 12454.5
$ws.Range("I116").Value = 12136
$ws.Range("K116").Value = 12136
$ws.Range("M116").Value = -8694
$ws.Range("H125").Value = 1928.7
$ws.Range("I125").Value = 704.3333
$ws.Range("K125").Value = 6338.9997
$ws.Range("M125").Value = -3878.9997
$ws.Range("H129").Value = 2220.35
$ws.Range("I129").Value = 1450.3
$ws.Range("J129").Value = 2990.4
$ws.Range("K129").Value = 4350.9
$ws.Range("L129").Value = 8971.200000000001
$ws.Range("M129").Value = 649.1000000000004
$ws.Range("N129").Value = -18971.2
$ws.Range("H132").Value = 1572
$ws.Range("I132").Value = 1424.3846
$ws.Range("K132").Value = 4273.1538
$ws.Range("M132").Value = -1743.1538
$ws.Range("H137").Value = 20002530
$ws.Range("J137").Value = 2809.1936
$ws.Range("L137").Value = 8427.5808
$ws.Range("N137").Value = -13527.5808
$ws.Range("H138").Value = 2938.1528
$ws.Range("I138").Value = 1519.4
$ws.Range("J138").Value = 3951.5476
$ws.Range("K138").Value = 4558.200000000001
$ws.Range("L138").Value = 11854.6428
$ws.Range("M138").Value = 581.7999999999993
$ws.Range("N138").Value = -22134.6428

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 5850042.5
$ws.Range("I74").Value = 6411676.5
$ws.Range("K74").Value = 6411676.5
$ws.Range("M74").Value = -6410802.5
$ws.Range("H77").Value = 5850042.5
$ws.Range("I77").Value = 6411676.5
$ws.Range("K77").Value = 32058382.5
$ws.Range("M77").Value = -32054014.5
$ws.Range("H132").Value = 3394.283
$ws.Range("I132").Value = 3004.0444
$ws.Range("K132").Value = 9012.1332
$ws.Range("M132").Value = -6482.1332

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 9991.5
$ws.Range("I82").Value = 9991.5
$ws.Range("K82").Value = 9991.5
$ws.Range("M82").Value = -9608.5
$ws.Range("H85").Value = 9991.5
$ws.Range("I85").Value = 9991.5
$ws.Range("K85").Value = 9991.5
$ws.Range("M85").Value = -8665.5
$ws.Range("H88").Value = 348458.56
$ws.Range("J88").Value = 348458.56
$ws.Range("L88").Value = 348458.56
$ws.Range("N88").Value = -349270.56
$ws.Range("H91").Value = 348458.56
$ws.Range("J91").Value = 348458.56
$ws.Range("L91").Value = 348458.56
$ws.Range("N91").Value = -351266.56
$ws.Range("H134").Value = 1497.0952
$ws.Range("I134").Value = 1055.5294
$ws.Range("K134").Value = 3166.5882
$ws.Range("M134").Value = -631.5881999999997

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 34662.57
$ws.Range("I31").Value = 3381.8333
$ws.Range("K31").Value = 3381.8333
$ws.Range("M31").Value = -3086.8333
$ws.Range("H34").Value = 34662.57
$ws.Range("I34").Value = 3381.8333
$ws.Range("K34").Value = 3381.8333
$ws.Range("M34").Value = -3179.8333
$ws.Range("H58").Value = 4357.7085
$ws.Range("I58").Value = 2027.1765
$ws.Range("K58").Value = 2027.1765
$ws.Range("M58").Value = -1824.1765
$ws.Range("H99").Value = 3976.111
$ws.Range("I99").Value = 3973.125
$ws.Range("K99").Value = 3973.125
$ws.Range("M99").Value = -2475.125
$ws.Range("H107").Value = 1934.125
$ws.Range("I107").Value = 1697.1538
$ws.Range("J107").Value = 2961
$ws.Range("K107").Value = 1697.1538
$ws.Range("L107").Value = 2961
$ws.Range("M107").Value = 222.8462
$ws.Range("N107").Value = -6801
$ws.Range("H126").Value = 3976.111
$ws.Range("I126").Value = 3973.125
$ws.Range("K126").Value = 11919.375
$ws.Range("M126").Value = -9449.375
$ws.Range("H132").Value = 2009.8948
$ws.Range("I132").Value = 1660.9706
$ws.Range("K132").Value = 4982.9118
$ws.Range("M132").Value = -2452.9118
$ws.Range("H134").Value = 7749.921
$ws.Range("I134").Value = 7288.3335
$ws.Range("K134").Value = 21865.0005
$ws.Range("M134").Value = -19330.0005
$ws.Range("H136").Value = 4357.7085
$ws.Range("I136").Value = 2027.1765
$ws.Range("K136").Value = 6081.529500000001
$ws.Range("M136").Value = -3531.529500000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 8028529
$ws.Range("I4").Value = 6729725
$ws.Range("K4").Value = 20189175
$ws.Range("M4").Value = -20189063
$ws.Range("H5").Value = 2381.6667
$ws.Range("I5").Value = 1263.6154
$ws.Range("J5").Value = 4198.5
$ws.Range("K5").Value = 3790.8462
$ws.Range("L5").Value = 12595.5
$ws.Range("M5").Value = -3678.8462
$ws.Range("N5").Value = -12819.5
$ws.Range("H102").Value = 6214.5
$ws.Range("J102").Value = 6214.5
$ws.Range("L102").Value = 18643.5
$ws.Range("N102").Value = -23511.5
$ws.Range("H122").Value = 10532679
$ws.Range("J122").Value = 12989036
$ws.Range("L122").Value = 116901324
$ws.Range("N122").Value = -116906224
$ws.Range("H135").Value = 2381.6667
$ws.Range("I135").Value = 1263.6154
$ws.Range("J135").Value = 4198.5
$ws.Range("K135").Value = 11372.5386
$ws.Range("L135").Value = 37786.5
$ws.Range("M135").Value = -8837.5386
$ws.Range("N135").Value = -42856.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1353.6072
$ws.Range("I97").Value = 1294.5652
$ws.Range("K97").Value = 1294.5652
$ws.Range("M97").Value = -798.5652
$ws.Range("H126").Value = 5055
$ws.Range("I126").Value = 3578.4
$ws.Range("J126").Value = 6285.5
$ws.Range("K126").Value = 10735.2
$ws.Range("L126").Value = 18856.5
$ws.Range("M126").Value = -8265.200000000001
$ws.Range("N126").Value = -23796.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3951.2173
$ws.Range("I68").Value = 3344.05
$ws.Range("K68").Value = 3344.05
$ws.Range("M68").Value = -2595.05
$ws.Range("H71").Value = 3951.2173
$ws.Range("I71").Value = 3344.05
$ws.Range("K71").Value = 16720.25
$ws.Range("M71").Value = -12976.25
$ws.Range("H116").Value = 289666.66
$ws.Range("J116").Value = 289666.66
$ws.Range("L116").Value = 289666.66
$ws.Range("N116").Value = -298844.66
$ws.Range("H132").Value = 7251.1904
$ws.Range("I132").Value = 6080.375
$ws.Range("K132").Value = 18241.125
$ws.Range("M132").Value = -15711.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 2000000
$ws.Range("J8").Value = 2000000
$ws.Range("L8").Value = 2000000
$ws.Range("N8").Value = -2000280
$ws.Range("H100").Value = 362.08334
$ws.Range("I100").Value = 363.85715
$ws.Range("J100").Value = 359.6
$ws.Range("K100").Value = 727.7143
$ws.Range("L100").Value = 719.2
$ws.Range("M100").Value = -186.7143
$ws.Range("N100").Value = -1801.2
$ws.Range("H132").Value = 3311.1562
$ws.Range("I132").Value = 1690.9
$ws.Range("K132").Value = 5072.700000000001
$ws.Range("M132").Value = -2542.700000000001
$ws.Range("H136").Value = 2293.024
$ws.Range("I136").Value = 1494.7028
$ws.Range("K136").Value = 4484.1084
$ws.Range("M136").Value = -1934.1084
